$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for rows 2-9, columns B,C,D,E,F,H (G is unchanged)
$data = @{
    2 = @{ B = 1300.617065429688; C = 0.9757;              D = 0.9520999789237976; E = 1.415799975395203;  F = 0.8181999921798706;  H = 1.0377 }
    3 = @{ B = 1262.348510742188; C = 1.0099;               D = 0.979;              E = 1.610100030899048;  F = 0.6743000149726868;  H = 1.2757 }
    4 = @{ B = 829.6500244140625; C = 0.9692;               D = 0.9272;             E = 1.661700010299683;  F = 0.7473000288009644;  H = 0.8174 }
    5 = @{ B = 808.3590087890625; C = 0.8527;               D = 0.8549;             E = 1.171900033950806;  F = 0.5960000157356262;  H = 0.1764 }
    6 = @{ B = 1085.690307617188; C = 0.8589;               D = 0.8578;             E = 1.06879997253418;   F = 0.6510000228881836;  H = 0.2027 }
    7 = @{ B = 835.6405029296875; C = 0.8415;               D = 0.8392999768257141; E = 1.060899972915649;  F = 0.675000011920929;   H = 0.0383 }
    8 = @{ B = 909.531005859375;  C = 0.8149999999999999;   D = 0.8132;             E = 1.080299973487854;  F = 0.6909999847412109;  H = -0.1924 }
    9 = @{ B = 7031.8369140625;   C = 0.9063;                D = 0.8848;             E = 1.661700010299683;  F = 0.5960000157356262;  H = 3.3558 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("H$row").Value = $vals.H
}
